$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New summary table (columns G/H/I) -------------------------------------
# Insertion order matters for shared-string dedup indices: the labels in
# column G are entered first (rows 2-4), then the column-H header, matching
# how the authoring session built up the shared-strings table.
$ws.Range("G2").Value = "1984-2004"
$ws.Range("G3").Value = "2005-2007"
$ws.Range("G4").Value = "2008-2020"
$ws.Range("H1").Value = "km2 mean"

$ws.Range("H2").Formula = "=AVERAGE(D2:D8)"
$ws.Range("H3").Formula = "=AVERAGE(D9:D11)"
$ws.Range("H4").Formula = "=AVERAGE(D12:D24)"

# Percent-change cell, formatted as a percentage (numFmtId 10 -> 0.00%)
$ws.Range("I4").Formula = "=(H4-H2)/H2"
$ws.Range("I4").NumberFormat = "0.00%"

# --- Selection / view -------------------------------------------------------
$ws.Range("E9:E21").Select() | Out-Null

# --- Move/resize the existing chart to make room for the new table ---------
$co = $ws.ChartObjects().Item(1)
$co.Left = 486.4022650098425
$co.Top = 153
$co.Width = 802.375
$co.Height = 216
